$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9019787907600403
$ws.Range("B1").Value = 1.64541482925415
$ws.Range("C1").Value = 4.270478248596191
$ws.Range("D1").Value = 2.745417356491089
$ws.Range("E1").Value = 0.6659830808639526
